# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 8
    3  = 3
    4  = 4
    5  = 10
    6  = 9
    7  = 3
    8  = 9
    9  = 3
    10 = 5
    11 = 2
    12 = 7
    13 = 7
    14 = 7
    15 = 5
    16 = 6
    17 = 9
    18 = 14
    19 = 3
    20 = 5
    21 = 4
    22 = 6
    23 = 6
    24 = 5
    25 = 3
    26 = 10
    27 = 8
    28 = 8
    29 = 6
    30 = 7
    31 = 7
    32 = 2
    33 = 7
    34 = 6
    35 = 3
    36 = 3
    37 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
